$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 653.9476226199286
$ws.Range("C2").Value = 136.3288556572895
$ws.Range("D2").Value = 85.96170863986495

$ws.Range("B3").Value = 653.9476226199286
$ws.Range("C3").Value = 136.3288556572895
$ws.Range("D3").Value = 85.96170863986495

$ws.Range("B4").Value = 28.80387728896227
$ws.Range("C4").Value = 46.71042583927353
$ws.Range("D4").Value = 0.826043741803197

$ws.Range("B5").Value = 625.1437453309662
$ws.Range("C5").Value = 89.61842981801597
$ws.Range("D5").Value = 85.13566489806176

$ws.Range("B6").Value = 682.7514999088908
$ws.Range("C6").Value = 183.039281496563
$ws.Range("D6").Value = 86.78775238166816
